$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.453.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.672.83"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "651.57"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -4.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.20"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.04"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.438"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000230"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.293.57"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.44"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.677.58"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.429.07"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.95"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.41"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "465.18"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.67"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.12%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.46"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.820.60"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.81"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.87"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.62"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.85%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.98"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.58"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.43"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.664.36"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.37"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.94"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "179.11"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.10%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.19"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0893"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.90"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.71"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.06"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000265"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -5.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.78"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.05"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.19%  "
